# ---------------------------------------------------------------------------
# Template edit: "added 34 pages to template"
#   1. Re-balance the two right-hand columns of the header table
#      (3774/3999 dxa -> 3773/4000 dxa).
#   2. Tidy the "Закінчено" run (was split across 3 runs, now a single run
#      with the same visible text).
#   3. Bump the page count "На  33  аркушах" -> "На  34  аркушах".
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Column widths -------------------------------------------------------
# Row 3 of table 1 holds the 1868 / 3774 / 3999 dxa columns; only the middle
# and right ones move (188.7pt -> 188.65pt, 199.95pt -> 200pt).
$headerTable = $d.Tables.Item(1)
$dataRow = $headerTable.Rows.Item(3)
$dataRow.Cells.Item(2).Width = 3773 / 20.0
$dataRow.Cells.Item(3).Width = 4000 / 20.0

# --- 2. Merge the "Закінчено" runs -----------------------------------------
# The three runs already spell out "Закінчено:   5 березня 2040 р." when
# concatenated, so there is nothing to retype -- just force Word to
# collapse them into one run by touching the tail of the text (inserting
# then deleting a throw-away character keeps the visible text identical).
$find = $d.Content.Find
$find.Execute("Закінчено:   5 березня 2040 р.", $false, $false, $false, `
               $false, $false, $true, 1, $false, "", 0)
if ($find.Found) {
    $finished = $find.Parent
    $tail = $finished.End
    $d.Range($tail, $tail).InsertAfter("~")
    $d.Range($tail, $tail + 1).Delete()
}

# --- 3. Bump the page count 33 -> 34 ----------------------------------------
$find2 = $d.Content.Find
$find2.Execute("На  33  аркушах", $false, $false, $false, $false, $false, `
                $true, 1, $false, "", 0)
if ($find2.Found) {
    $pages = $find2.Parent
    $offset = $pages.Text.IndexOf("33")
    $digits = $d.Range($pages.Start + $offset, $pages.Start + $offset + 2)
    $digits.Text = "34"
}
